$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.050.02'
$ws.Range("E2").Value = '  +1.15%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.621.68'
$ws.Range("E3").Value = '  +1.19%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.58'
$ws.Range("E5").Value = '  +3.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.28'
$ws.Range("E6").Value = '  +0.24%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  +1.83%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.70'
$ws.Range("E9").Value = '  +0.71%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +4.02%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.348'
$ws.Range("E11").Value = '  +0.53%  '

# Row 12
$ws.Range("E12").Value = '  +0.39%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.080.21'
$ws.Range("E13").Value = '  +1.10%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '61.050.62'
$ws.Range("E14").Value = '  +1.08%  '

# Row 15
$ws.Range("E15").Value = '  +0.26%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000143'
$ws.Range("E16").Value = '  +1.76%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.626.96'
$ws.Range("E17").Value = '  +0.98%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.79'
$ws.Range("E18").Value = '  +0.61%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '355.80'
$ws.Range("E19").Value = '  +2.15%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.61'
$ws.Range("E20").Value = '  +0.94%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("E21").Value = '  +2.16%  '

# Row 22
$ws.Range("E22").Value = '  +0.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.69'
$ws.Range("E23").Value = '  +2.00%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.430'
$ws.Range("E24").Value = '  +2.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  +1.29%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.735.77'
$ws.Range("E26").Value = '  +1.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.995'
$ws.Range("E27").Value = '  -0.40%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0862'
$ws.Range("E28").Value = '  +1.88%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("E29").Value = '  +0.12%  '

# Row 30
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.22'
$ws.Range("E31").Value = '  +8.53%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.48'
$ws.Range("E32").Value = '  +0.60%  '

# Row 33
$ws.Range("E33").Value = '  +3.19%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.63'
$ws.Range("E34").Value = '  -1.05%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.17'
$ws.Range("E35").Value = '  +3.18%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.21'
$ws.Range("E36").Value = '  +1.48%  '

# Row 37
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.935'
$ws.Range("E37").Value = '  +10.36%  '

# Row 38
$ws.Range("B38").Value = 'SuiNetwork'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.902'
$ws.Range("E38").Value = '  +4.87%  '

# Row 39
$ws.Range("E39").Value = '  +1.44%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.80'
$ws.Range("E40").Value = '  +0.65%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '294.41'
$ws.Range("E41").Value = '  -1.40%  '

# Row 42
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.639'
$ws.Range("E42").Value = '  +2.81%  '

# Row 43
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.103'
$ws.Range("E43").Value = '  +2.70%  '

# Row 44
$ws.Range("E44").Value = '  +0.43%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.01%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.77'
$ws.Range("E46").Value = '  -0.25%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.92'
$ws.Range("E47").Value = '  +2.00%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0239'
$ws.Range("E48").Value = '  +2.21%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.27'
$ws.Range("E49").Value = '  +7.10%  '

# Row 50
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.35'
$ws.Range("E50").Value = '  +0.51%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.981.80'
$ws.Range("E51").Value = '  -0.70%  '
